$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" date column (C) for rows 2-9 from 45174 (2023-09-05)
# to 45175 (2023-09-06), per the automatic update of files.
foreach ($row in 2..9) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45174) {
        $cell.Value2 = 45175
    }
}
